$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Nov 01 15:54:00 EDT 2023"
$ws.Range("B3").Value = "Wed Nov 01 15:54:13 EDT 2023"
$ws.Range("B4").Value = "Wed Nov 01 15:54:26 EDT 2023"
